$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 14. Existing rows 14-24 (the
# weekly price records) shift down to become rows 15-25; the new row 14
# inherits the row-13 formatting (keeps the date column's style).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with this week's record.
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C14").Value = "Los Lagos"
$ws.Range("D14").Value = 44827
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 100112035
$ws.Range("G14").Value = "Bruselas (repollito)"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 90
$ws.Range("K14").Value = 22000
$ws.Range("L14").Value = 22000
$ws.Range("M14").Value = 22000
$ws.Range("N14").Value = "$/malla 15 kilos"
$ws.Range("O14").Value = "Provincia de Quillota"
$ws.Range("P14").Value = 1467
$ws.Range("Q14").Value = 15
$ws.Range("R14").Value = "Hortaliza"
